$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 37 (shift cells down), pushing old rows 37:84 to 38:85
$ws.Range("A37:R37").Insert(-4121)

# Copy formatting (including date style) from the row below into the new row 37
$ws.Range("A38:R38").Copy()
$ws.Range("A37:R37").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row 37 with the provided values
$ws.Cells.Item(37, 1).Value = 4
$ws.Cells.Item(37, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(37, 3).Value = "Los Lagos"
$ws.Cells.Item(37, 4).Value = 44638
$ws.Cells.Item(37, 5).Value = 10
$ws.Cells.Item(37, 6).Value = 100112026
$ws.Cells.Item(37, 7).Value = "Haba"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 50
$ws.Cells.Item(37, 11).Value = 26000
$ws.Cells.Item(37, 12).Value = 26000
$ws.Cells.Item(37, 13).Value = 26000
$ws.Cells.Item(37, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(37, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(37, 16).Value = 1040
$ws.Cells.Item(37, 17).Value = 25
$ws.Cells.Item(37, 18).Value = "Hortaliza"
